$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Tommy" data points (column L/M) mirroring the existing
# Richi (F/G) / Laura (I/J) rows that were already present.
$ws.Range("L6").Value = "34s"
$ws.Range("M6").Value = 97598

$ws.Range("L15").Value = "36s"
$ws.Range("M15").Value = 82832

$ws.Range("L28").Value = "1min49"
$ws.Range("M28").Value = 90912

# Update the saved view state (selection / scroll position).
$ws.Range("L29").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
